$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Select Menu" widget section — inserted after the existing "Menu" section
# (rows 170-177), following the same layout as the rest of the sheet:
#   column A = ElementID, column B = ElementPath, column C = Method (By.xpath)

# Section header row (ElementID only)
$ws.Range("A170").Value = "selectMenu"

$rows = @(
    @{ Row = 171; Id = "selectMenuNav";          Path = '//div[@class=''element-list collapse show'']//li[@id=''item-8'']' },
    @{ Row = 172; Id = "selectMenuElement";      Path = '//*[@id="withOptGroup"]/div/div[2]/div' },
    @{ Row = 173; Id = "selectedOptionElement";  Path = '//*[@id="react-select-23-option-2"]' },
    @{ Row = 174; Id = "selectMenuScroll";       Path = "//h1[normalize-space()='Select Menu']" },
    @{ Row = 175; Id = "selectOne";              Path = '//*[@id="selectOne"]/div[1]/div[2]/div' },
    @{ Row = 176; Id = "OldStyleMenuSelect";     Path = '//*[@id="oldSelectMenu"]' },
    @{ Row = 177; Id = "MultiMenuSelect";        Path = '//*[@id="selectMenuContainer"]/div[7]/div/div/div' }
)

foreach ($r in $rows) {
    $ws.Range("A$($r.Row)").Value = $r.Id
    $ws.Range("B$($r.Row)").Value = $r.Path
    $ws.Range("C$($r.Row)").Value = "By.xpath"
    $ws.Range("C$($r.Row)").Font.Name = "Calibri"
}

# Update the sheet's active selection to match the new bottom of the table
$ws.Range("B176").Select()
